# Weekly update: insert 2 new data rows (new week's prices) above the
# existing row 592, shifting the old rows 592-672 down to 594-674.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 592 (old rows 592-672 shift to 594-674).
$ws.Range("A592:R593").EntireRow.Insert()

# --- New row 592 ---
$ws.Cells.Item(592, 1).Value = 3
$ws.Cells.Item(592, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(592, 3).Value = "Coquimbo"
$ws.Cells.Item(592, 4).Value = 44984
$ws.Cells.Item(592, 5).Value = 5
$ws.Cells.Item(592, 6).Value = 100112032
$ws.Cells.Item(592, 7).Value = "Zapallo italiano"
$ws.Cells.Item(592, 8).Value = "Sin especificar"
$ws.Cells.Item(592, 9).Value = "Primera"
$ws.Cells.Item(592, 10).Value = 170
$ws.Cells.Item(592, 11).Value = 4000
$ws.Cells.Item(592, 12).Value = 4500
$ws.Cells.Item(592, 13).Value = 4235
$ws.Cells.Item(592, 14).Value = "$/caja 36 unidades"
$ws.Cells.Item(592, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(592, 16).Value = 118
$ws.Cells.Item(592, 17).Value = 36
$ws.Cells.Item(592, 18).Value = "Hortaliza"

# --- New row 593 ---
$ws.Cells.Item(593, 1).Value = 3
$ws.Cells.Item(593, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(593, 3).Value = "Coquimbo"
$ws.Cells.Item(593, 4).Value = 44984
$ws.Cells.Item(593, 5).Value = 5
$ws.Cells.Item(593, 6).Value = 100112032
$ws.Cells.Item(593, 7).Value = "Zapallo italiano"
$ws.Cells.Item(593, 8).Value = "Sin especificar"
$ws.Cells.Item(593, 9).Value = "Primera"
$ws.Cells.Item(593, 10).Value = 160
$ws.Cells.Item(593, 11).Value = 6600
$ws.Cells.Item(593, 12).Value = 7000
$ws.Cells.Item(593, 13).Value = 6812
$ws.Cells.Item(593, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(593, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(593, 16).Value = 114
$ws.Cells.Item(593, 17).Value = 60
$ws.Cells.Item(593, 18).Value = "Hortaliza"
